# se cambia data para smoke en qa
$wb = $excel.ActiveWorkbook

$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsHogar  = $wb.Worksheets.Item("DatosHogar")
$wsMotor  = $wb.Worksheets.Item("DatosMotor")
$wsAP     = $wb.Worksheets.Item("DatosAP")

# DatosCuenta: new account row data
$wsCuenta.Range("A2").Value = "SmokeQAUno"
$wsCuenta.Range("B2").Value = "SmokeQANameUno"
$wsCuenta.Range("C2").Value = 27100114
$wsCuenta.Range("D2").Value = 116
$wsCuenta.Activate()
$wsCuenta.Range("D2").Select()

# DatosHogar: increment NvoNro
$wsHogar.Range("A2").Value = 635

# DatosMotor: new plate/motor/chassis data
$wsMotor.Range("A2").Value = "SMA017"
$wsMotor.Range("B2").Value = "ABC12SSMA017"
$wsMotor.Range("C2").Value = "ZAZ123SSMA017"

# DatosAP: new document number
$wsAP.Range("A2").Value = 21200117

# restore original active sheet/tab (DatosAP) and its selection
$wsAP.Activate()
$wsAP.Range("A3").Select()
